$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("7:8").Insert()
$ws.Range("A7").Value = "CHE"
$ws.Range("B7").Value = "conv_chp_biogas"
$ws.Range("C7").Value = "input"
$ws.Range("F7").Value = "biogas"
$ws.Range("G7").Value = 1

$ws.Range("A8").Value = "CHE"
$ws.Range("B8").Value = "conv_chp_biogas"
$ws.Range("C8").Value = "output"
$ws.Range("F8").Value = "elecsupply"
$ws.Range("G8").Value = 1

$ws.Range("D7").Value = "configuration_fxe"
$ws.Range("D8").Value = "configuration_fxe"

$ws.Range("G9").Select()
$ws.Columns.Item(4).EntireColumn.AutoFit()
